# Write the 3 header cells (A1:C1) and apply the "table header" formatting:
# bold font, thin border on all sides, centered horizontally, top-aligned
# vertically — matching the sample_data.xlsx table-template header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("col1", "col2", "col3")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous -> thin border
}
